# "revising decision models data"
#
# This script reproduces, via the Excel COM object model, the content
# changes made to class_materials/data/decision_model.xlsx:
#   - fix the "Agents Avaialble" typo -> "Agents Available"
#   - replace the terse "Min." / "Max." column headers (F1/G1) with
#     full descriptive labels, matching the wrap-text style already
#     used by the C1 header
#   - grow row 1 slightly to better fit the new, longer header text
#   - add a new, currency-formatted (empty) cell at I9
#   - move the active selection to I9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Agents Avaialble" -> "Agents Available" typo (used by G8).
$ws.Range("G8").Value = "Agents Available"

# Replace the short "Min." / "Max." headers with full descriptive text.
$ws.Range("F1").Value = "Minimum Number of Agents Assigned"
$ws.Range("G1").Value = "Maximum Number of Agents Assigned"

# Give the new, longer headers the same wrapped look as the other
# multi-line header (C1), and grow the header row to fit.
$ws.Range("F1:G1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 63

# Add a new, empty, currency-formatted cell below the Revenue label.
$ws.Range("I9").Style = "Currency"

# Leave the active selection on the newly added cell.
$null = $ws.Range("I9").Select()
